# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet named "2022-Q1" right before the "总计" sheet,
#    formatted like the other quarterly sheets (2020-Q4 / 2021-Q2 / 2021-Q3),
#    and fill it with the Q1-2022 fund holdings table.
# 2. Update the "总计" (grand total) summary sheet with a new leading row
#    for 2022-Q1 (pushing the older quarters down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet positioned right before "总计"
# ---------------------------------------------------------------------

$zongji = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($zongji)
$q1.Name = "2022-Q1"

# Recreate the sheetPr/outlinePr defaults used by every other sheet in the
# workbook (summaryBelow / summaryRight).
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = -4152

# Borrow the header/index-column look (bold, centered, thin border) from the
# "2021-Q3" sheet so the new tab matches its siblings pixel-for-pixel,
# instead of re-building the style by hand.
$ref = $wb.Worksheets.Item("2021-Q3")
$ref.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$ref.Range("A2").Copy()
$q1.Range("A2:A8").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B:G hold text (fund codes / decimal strings that must keep
# trailing & leading zeros), so force a text number format before writing
# them — otherwise COM auto-coerces "007950" -> 7950 and "181.00" -> 181.
$q1.Range("B2:G8").NumberFormat = "@"

$rows = @(
    @(0, "510880", "华泰柏瑞上证红利ETF",        "181.00", "97.22", "2.27", "4.1087", 9),
    @(1, "009931", "淳厚欣享一年持有期混合A",      "10.95",  "88.81", "2.68", "0.2935", 6),
    @(2, "008186", "淳厚信睿核心精选混合A",        "8.49",   "88.47", "2.79", "0.2369", 6),
    @(3, "009939", "淳厚欣享一年持有期混合C",      "1.96",   "88.81", "2.68", "0.0525", 6),
    @(4, "008187", "淳厚信睿核心精选混合C",        "1.65",   "88.47", "2.79", "0.0460", 6),
    @(5, "001917", "招商量化精选股票A",            "2.33",   "94.20", "1.43", "0.0333", 5),
    @(6, "007950", "招商量化精选股票C",            "0.56",   "94.20", "1.43", "0.0080", 5)
)

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# The "@" number format above was only a write-time guard against numeric
# auto-coercion; the source table has no special formatting on these data
# cells, so drop back to the default "Normal" style now that the text is
# safely in place.
$q1.Range("B2:G8").Style = "Normal"

# ---------------------------------------------------------------------
# Step 2: update the "总计" sheet with a new 2022-Q1 row on top, shifting
# the existing quarters down by one row.
# ---------------------------------------------------------------------

$tot = $wb.Worksheets.Item("总计")

$tot.Cells.Item(2, 1).Value = 0
$tot.Cells.Item(2, 2).Value = "2022-Q1"
$tot.Cells.Item(2, 3).Value = 7
$tot.Cells.Item(2, 4).Value = 4.78

$tot.Cells.Item(3, 1).Value = 1
$tot.Cells.Item(3, 2).Value = "2021-Q3"
$tot.Cells.Item(3, 3).Value = 2
$tot.Cells.Item(3, 4).Value = 4.17

$tot.Cells.Item(4, 1).Value = 2
$tot.Cells.Item(4, 2).Value = "2021-Q2"
$tot.Cells.Item(4, 3).Value = 1
$tot.Cells.Item(4, 4).Value = 0.29

$tot.Cells.Item(5, 1).Value = 3
$tot.Cells.Item(5, 2).Value = "2020-Q4"
$tot.Cells.Item(5, 3).Value = 6
$tot.Cells.Item(5, 4).Value = 4.19

# Row 5 is brand new territory for this sheet - copy the index-column style
# (bold/centered/bordered) from A2 so it matches A2:A4.
$tot.Range("A2").Copy()
$tot.Range("A5").PasteSpecial(-4122)
